# FantasyTeamPoints.xlsx update
# - renames Sheet1 -> Oct24Data
# - adds a new sheet Nov05Data (after Oct24Data) with the latest pull of
#   fantasy hockey stats
# - leaves Oct24Data's selection parked at the header row, Nov05Data
#   becomes the active/selected tab

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Oct24Data"

# reset the old sheet's selection to the header row (it's no longer the
# tab the user left open)
$ws1.Range("A1:K1").Select()

# --- add the new sheet right after Oct24Data -----------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Nov05Data"

# --- header row -----------------------------------------------------------------
$headers = @("name", "position_1", "position_2", "team", "games_7", "points_7", "games_14", "points_14", "games_30", "points_30", "games_this_week")
for ($col = 1; $col -le $headers.Count; $col++) {
    $ws2.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- data rows --------------------------------------------------------------
# name, position_1, position_2, team, games_7, points_7, games_14, points_14, games_30, points_30, games_this_week
$data = @(
    @("Pettersson", "C", $null, "VAN", 3, 61, 4, 61, 8, 111.8, 3),
    @("Monahan", "C", $null, "CGY", 4, 64.6, 7, 78.7, 14, 139.2, 3),
    @("Skinner", "C", "L", "BUF", 4, 46.1, 6, 80.6, 14, 163.2, 2),
    @("Lee", "L", $null, "NYI", 3, 29.3, 6, 87.5, 12, 135.3, 3),
    @("Gourde", "L", "R", "TB", 4, 51.5, 7, 62, 14, 124.6, 3),
    @("Wheeler", "C", "R", "WPG", 2, 37.6, 6, 79, 13, 119, 2),
    @("Ellis", "D", $null, "NSH", 3, 28.6, 6, 42.2, 13, 93.8, 2),
    @("Ekholm", "D", $null, "NSH", 3, 16.6, 6, 27, 13, 84.7, 2),
    @("Suter", "D", $null, "MIN", 3, 21.5, 5, 39, 12, 107.8, 4),
    @("Seabrook", "D", $null, "CHI", 3, 17.3, 7, 52.9, 14, 89, 2),
    @("Marner", "C", "R", "TOR", 3, 26.7, 5, 46.1, 13, 141.6, 3),
    @("Kane", "R", $null, "CHI", 2, 8.4, 6, 58.9, 13, 151.2, 2),
    @("Hall", "L", $null, "NJ", 3, 10, 5, 38.5, 11, 105.4, 4),
    @("Backstrom", "C", $null, "WAS", 2, 17.4, 5, 34.1, 10, 84.2, 4),
    @("Hellebuyck", "G", $null, "WPG", 2, 21.8, 5, 49.6, 10, 92, 2),
    @("Bishop", "G", $null, "DAL", 2, 35.8, 5, 68.6, 9, 102.4, 4)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws2.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws2.Cells.Item($rowIndex, 2).Value = $row[1]
    if ($row[2] -ne $null) {
        $ws2.Cells.Item($rowIndex, 3).Value = $row[2]
    }
    $ws2.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws2.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws2.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws2.Cells.Item($rowIndex, 6).NumberFormat = "0.00"
    $ws2.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws2.Cells.Item($rowIndex, 8).Value = $row[7]
    $ws2.Cells.Item($rowIndex, 8).NumberFormat = "0.00"
    $ws2.Cells.Item($rowIndex, 9).Value = $row[8]
    $ws2.Cells.Item($rowIndex, 10).Value = $row[9]
    $ws2.Cells.Item($rowIndex, 10).NumberFormat = "0.00"
    $ws2.Cells.Item($rowIndex, 11).Value = $row[10]
    # column L mirrors the source workbook: present, numeric-formatted, empty
    $ws2.Cells.Item($rowIndex, 12).NumberFormat = "0.00"
    $rowIndex++
}

# --- selection / active tab ------------------------------------------------
$ws2.Range("J26").Select()
$ws2.Activate()
